$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 21) mirroring the existing daily COVID data table.

# Copy row 20's formatting down first so the new date cell reuses the same
# style (numFmtId 14, "m/d/yyyy") instead of Excel inventing a brand-new one.
$ws.Range("B20").Copy()
$ws.Range("B21").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A21").Value = "MNE"
$ws.Range("B21").Value = 43918
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 385
$ws.Range("E21").Value = 84
$ws.Range("F21").Formula = "=(E21-E20)/E20"
$ws.Range("G21").Formula = "=E21-E20"
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = 6278

# Mirror the selection state recorded after the edit
$ws.Range("D22").Select()
